$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 476; existing rows 476:519 shift down to 477:520
$ws.Rows("476:476").Insert()

# Populate the newly inserted row 476 with the new record
$ws.Range("A476").Value = 5
$ws.Range("B476").Value = "Macroferia Regional de Talca"
$ws.Range("C476").Value = "Maule"
$ws.Range("D476").Value = 45132
$ws.Range("E476").Value = 7
$ws.Range("F476").Value = 100112003
$ws.Range("G476").Value = "Ajo"
$ws.Range("H476").Value = "Chino"
$ws.Range("I476").Value = "Primera"
$ws.Range("J476").Value = 300
$ws.Range("K476").Value = 20000
$ws.Range("L476").Value = 20000
$ws.Range("M476").Value = 20000
$ws.Range("N476").Value = "`$/malla 10 kilos"
$ws.Range("O476").Value = "China"
$ws.Range("P476").Value = 2000
$ws.Range("Q476").Value = 10
$ws.Range("R476").Value = "Hortaliza"
